{"js": "// Revert \"More final work\": remove the bullet content that was added to\n// Holly Bennett's \"Contribution\" cell, restoring the cell to its prior\n// (effectively empty) state while leaving the trailing blank paragraphs.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No tables found in the document.\");\n}\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Locate the row whose first cell contains \"Holly Bennett\".\nlet targetRow = null;\nfor (let r = 0; r < rows.items.length; r++) {\n  const row = rows.items[r];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  const nameCell = cells.items[0];\n  nameCell.body.load(\"text\");\n  await context.sync();\n\n  if (nameCell.body.text.trim() === \"Holly Bennett\") {\n    targetRow = row;\n    break;\n  }\n}\n\nif (!targetRow) {\n  throw new Error('Could not find the \"Holly Bennett\" row.');\n}\n\nconst rowCells = targetRow.cells;\nrowCells.load(\"items\");\nawait context.sync();\n\n// The contribution text lives in the second cell of the row.\nconst contributionCell = rowCells.items[1];\nconst paragraphs = contributionCell.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Delete every paragraph that still has text or list-paragraph formatting\n// applied to it (the four bullet items plus the trailing empty\n// ListParagraph), leaving only the plain blank paragraphs behind.\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const para = items[i];\n  if (para.text !== \"\") {\n    para.delete();\n  }\n}\nawait context.sync();\n\n// Also clean up any now-empty paragraph that still carries the\n// ListParagraph style (no visible text, but leftover numbering/style).\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst remaining = paragraphs.items;\nfor (let i = 0; i < remaining.length; i++) {\n  const para = remaining[i];\n  if (para.text === \"\" && para.style === \"List Paragraph\") {\n    para.delete();\n  }\n}\nawait context.sync();\n", "ps1": "# Revert \"More final work\": remove the bullet content that was added to\n# Holly Bennett's \"Contribution\" cell, restoring the cell to its prior\n# (effectively empty) state while leaving the trailing blank paragraphs.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Locate the row whose first cell is \"Holly Bennett\".\n$targetRowIndex = -1\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $cellText = ($t.Cell($r, 1).Range.Text -replace \"[\\r\\a\\x07]\", \"\").Trim()\n    if ($cellText -eq \"Holly Bennett\") {\n        $targetRowIndex = $r\n        break\n    }\n}\n\nif ($targetRowIndex -eq -1) {\n    throw 'Could not find the \"Holly Bennett\" row.'\n}\n\n# Walk the whole-document paragraph collection (obtained from a fresh\n# Range() \u2014 going through a Table/Cell object first and then reusing\n# $d.Paragraphs leaves stale per-item Range data in this host) to find the\n# \"Holly Bennett\" paragraph itself, then collect the contiguous run of\n# \"List Paragraph\" styled paragraphs that follow it. That run is the four\n# bullet items plus the trailing empty ListParagraph that need to go.\n$paras = $d.Range().Paragraphs\n$count = $paras.Count\n\n$nameIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $paras.Item($i)\n    $txt = ($p.Range.Text -replace \"[\\r\\a\\x07]\", \"\").Trim()\n    if ($txt -eq \"Holly Bennett\") {\n        $nameIdx = $i\n        break\n    }\n}\n\nif ($nameIdx -eq -1) {\n    throw 'Could not find the \"Holly Bennett\" paragraph.'\n}\n\n$blockEnd = $nameIdx\nfor ($i = $nameIdx + 1; $i -le $count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Style.NameLocal -eq \"List Paragraph\") {\n        $blockEnd = $i\n    }\n    else {\n        break\n    }\n}\n\n# Delete from last to first so the earlier indices stay valid.\nfor ($i = $blockEnd; $i -ge ($nameIdx + 1); $i--) {\n    $paras.Item($i).Range.Delete()\n}\n"}
